$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 40.91514966666667
$ws.Range("H2").Value2 = 122.745449
$ws.Range("I2").Value2 = 0.8529192913871414
$ws.Range("J2").Value2 = 0.8529192913871415
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 115.5575153333333
$ws.Range("N2").Value2 = 346.672546
$ws.Range("O2").Value2 = 0.9048104954928987
$ws.Range("P2").Value2 = 0.9048104954928987
$ws.Range("Q2").Value2 = 4728.053034971463
$ws.Range("R2").Value2 = 42552.47731474316
$ws.Range("S2").Value2 = 0.7717303266554515
$ws.Range("T2").Value2 = 0.7717303266554516
$ws.Range("G3").Value2 = 40.91514966666667
$ws.Range("H3").Value2 = 122.745449
$ws.Range("I3").Value2 = 0.8529192913871414
$ws.Range("J3").Value2 = 0.8529192913871415
$ws.Range("M3").Value2 = 5.519651666666666
$ws.Range("O3").Value2 = 0.04321864090845719
$ws.Range("P3").Value2 = 0.04321864090845719
$ws.Range("Q3").Value2 = 225.8373740495328
$ws.Range("R3").Value2 = 2032.536366445795
$ws.Range("S3").Value2 = 0.03686201257835663
$ws.Range("T3").Value2 = 0.03686201257835663
$ws.Range("G4").Value2 = 40.91514966666667
$ws.Range("H4").Value2 = 122.745449
$ws.Range("I4").Value2 = 0.8529192913871414
$ws.Range("J4").Value2 = 0.8529192913871415
$ws.Range("M4").Value2 = 6.580297333333334
$ws.Range("N4").Value2 = 19.740892
$ws.Range("O4").Value2 = 0.05152345196666309
$ws.Range("P4").Value2 = 0.05152345196666309
$ws.Range("Q4").Value2 = 269.233850244501
$ws.Range("R4").Value2 = 2423.104652200509
$ws.Range("S4").Value2 = 0.0439453461412257
$ws.Range("T4").Value2 = 0.0439453461412257
$ws.Range("G5").Value2 = 40.91514966666667
$ws.Range("H5").Value2 = 122.745449
$ws.Range("I5").Value2 = 0.8529192913871414
$ws.Range("J5").Value2 = 0.8529192913871415
$ws.Range("K5").Value2 = 1
$ws.Range("L5").Value2 = 0.3333333333333333
$ws.Range("M5").Value2 = 0.057141
$ws.Range("N5").Value2 = 0.171423
$ws.Range("O5").Value2 = 0.0004474116319810314
$ws.Range("P5").Value2 = 0.0004474116319810314
$ws.Range("Q5").Value2 = 2.337932567103
$ws.Range("R5").Value2 = 21.041393103927
$ws.Range("S5").Value2 = 0.0003816060121076258
$ws.Range("T5").Value2 = 0.0003816060121076258
$ws.Range("I6").Value2 = 0.00344410114086962
$ws.Range("J6").Value2 = 0.003444101140869621
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 115.5575153333333
$ws.Range("N6").Value2 = 346.672546
$ws.Range("O6").Value2 = 0.9048104954928987
$ws.Range("P6").Value2 = 0.9048104954928987
$ws.Range("Q6").Value2 = 19.091950453312
$ws.Range("R6").Value2 = 171.827554079808
$ws.Range("S6").Value2 = 0.003116258859797899
$ws.Range("T6").Value2 = 0.0031162588597979
$ws.Range("I7").Value2 = 0.00344410114086962
$ws.Range("J7").Value2 = 0.003444101140869621
$ws.Range("M7").Value2 = 5.519651666666666
$ws.Range("O7").Value2 = 0.04321864090845719
$ws.Range("P7").Value2 = 0.04321864090845719
$ws.Range("Q7").Value2 = 0.9119347697599998
$ws.Range("R7").Value2 = 8.20741292784
$ws.Range("S7").Value2 = 0.0001488493704596519
$ws.Range("T7").Value2 = 0.0001488493704596519
$ws.Range("I8").Value2 = 0.00344410114086962
$ws.Range("J8").Value2 = 0.003444101140869621
$ws.Range("M8").Value2 = 6.580297333333334
$ws.Range("N8").Value2 = 19.740892
$ws.Range("O8").Value2 = 0.05152345196666309
$ws.Range("P8").Value2 = 0.05152345196666309
$ws.Range("Q8").Value2 = 1.087170404224
$ws.Range("R8").Value2 = 9.784533638016002
$ws.Range("S8").Value2 = 0.0001774519796999254
$ws.Range("T8").Value2 = 0.0001774519796999255
$ws.Range("I9").Value2 = 0.00344410114086962
$ws.Range("J9").Value2 = 0.003444101140869621
$ws.Range("K9").Value2 = 1
$ws.Range("L9").Value2 = 0.3333333333333333
$ws.Range("M9").Value2 = 0.057141
$ws.Range("N9").Value2 = 0.171423
$ws.Range("O9").Value2 = 0.0004474116319810314
$ws.Range("P9").Value2 = 0.0004474116319810314
$ws.Range("Q9").Value2 = 0.009440607456
$ws.Range("R9").Value2 = 0.084965467104
$ws.Range("S9").Value2 = 0.000001540930912144209
$ws.Range("T9").Value2 = 0.000001540930912144209
$ws.Range("E10").Value2 = 2
$ws.Range("F10").Value2 = 0.6666666666666666
$ws.Range("G10").Value2 = 0.4441646666666667
$ws.Range("H10").Value2 = 1.332494
$ws.Range("I10").Value2 = 0.009259079236881667
$ws.Range("J10").Value2 = 0.009259079236881667
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 115.5575153333333
$ws.Range("N10").Value2 = 346.672546
$ws.Range("O10").Value2 = 0.9048104954928987
$ws.Range("P10").Value2 = 0.9048104954928987
$ws.Range("Q10").Value2 = 51.32656527885823
$ws.Range("R10").Value2 = 461.9390875097241
$ws.Range("S10").Value2 = 0.008377712072130911
$ws.Range("T10").Value2 = 0.008377712072130911
$ws.Range("E11").Value2 = 2
$ws.Range("F11").Value2 = 0.6666666666666666
$ws.Range("G11").Value2 = 0.4441646666666667
$ws.Range("H11").Value2 = 1.332494
$ws.Range("I11").Value2 = 0.009259079236881667
$ws.Range("J11").Value2 = 0.009259079236881667
$ws.Range("M11").Value2 = 5.519651666666666
$ws.Range("O11").Value2 = 0.04321864090845719
$ws.Range("P11").Value2 = 0.04321864090845719
$ws.Range("Q11").Value2 = 2.451634242641111
$ws.Range("R11").Value2 = 22.06470818377
$ws.Range("S11").Value2 = 0.0004001648206817406
$ws.Range("T11").Value2 = 0.0004001648206817406
$ws.Range("E12").Value2 = 2
$ws.Range("F12").Value2 = 0.6666666666666666
$ws.Range("G12").Value2 = 0.4441646666666667
$ws.Range("H12").Value2 = 1.332494
$ws.Range("I12").Value2 = 0.009259079236881667
$ws.Range("J12").Value2 = 0.009259079236881667
$ws.Range("M12").Value2 = 6.580297333333334
$ws.Range("N12").Value2 = 19.740892
$ws.Range("O12").Value2 = 0.05152345196666309
$ws.Range("P12").Value2 = 0.05152345196666309
$ws.Range("Q12").Value2 = 2.922735571627556
$ws.Range("R12").Value2 = 26.30462014464801
$ws.Range("S12").Value2 = 0.0004770597243170001
$ws.Range("T12").Value2 = 0.0004770597243170001
$ws.Range("E13").Value2 = 2
$ws.Range("F13").Value2 = 0.6666666666666666
$ws.Range("G13").Value2 = 0.4441646666666667
$ws.Range("H13").Value2 = 1.332494
$ws.Range("I13").Value2 = 0.009259079236881667
$ws.Range("J13").Value2 = 0.009259079236881667
$ws.Range("K13").Value2 = 1
$ws.Range("L13").Value2 = 0.3333333333333333
$ws.Range("M13").Value2 = 0.057141
$ws.Range("N13").Value2 = 0.171423
$ws.Range("O13").Value2 = 0.0004474116319810314
$ws.Range("P13").Value2 = 0.0004474116319810314
$ws.Range("Q13").Value2 = 0.025380013218
$ws.Range("R13").Value2 = 0.228420118962
$ws.Range("S13").Value2 = 0.000004142619752014909
$ws.Range("T13").Value2 = 0.000004142619752014909
$ws.Range("G14").Value2 = 6.446186333333333
$ws.Range("H14").Value2 = 19.338559
$ws.Range("I14").Value2 = 0.1343775282351073
$ws.Range("J14").Value2 = 0.1343775282351073
$ws.Range("K14").Value2 = 3
$ws.Range("L14").Value2 = 1
$ws.Range("M14").Value2 = 115.5575153333333
$ws.Range("N14").Value2 = 346.672546
$ws.Range("O14").Value2 = 0.9048104954928987
$ws.Range("P14").Value2 = 0.9048104954928987
$ws.Range("Q14").Value2 = 744.9052760556905
$ws.Range("R14").Value2 = 6704.147484501214
$ws.Range("S14").Value2 = 0.1215861979055184
$ws.Range("T14").Value2 = 0.1215861979055184
$ws.Range("G15").Value2 = 6.446186333333333
$ws.Range("H15").Value2 = 19.338559
$ws.Range("I15").Value2 = 0.1343775282351073
$ws.Range("J15").Value2 = 0.1343775282351073
$ws.Range("M15").Value2 = 5.519651666666666
$ws.Range("O15").Value2 = 0.04321864090845719
$ws.Range("P15").Value2 = 0.04321864090845719
$ws.Range("Q15").Value2 = 35.58070313842722
$ws.Range("R15").Value2 = 320.226328245845
$ws.Range("S15").Value2 = 0.00580761413895917
$ws.Range("T15").Value2 = 0.005807614138959171
$ws.Range("G16").Value2 = 6.446186333333333
$ws.Range("H16").Value2 = 19.338559
$ws.Range("I16").Value2 = 0.1343775282351073
$ws.Range("J16").Value2 = 0.1343775282351073
$ws.Range("M16").Value2 = 6.580297333333334
$ws.Range("N16").Value2 = 19.740892
$ws.Range("O16").Value2 = 0.05152345196666309
$ws.Range("P16").Value2 = 0.05152345196666309
$ws.Range("Q16").Value2 = 42.41782273940311
$ws.Range("R16").Value2 = 381.760404654628
$ws.Range("S16").Value2 = 0.006923594121420463
$ws.Range("T16").Value2 = 0.006923594121420465
$ws.Range("G17").Value2 = 6.446186333333333
$ws.Range("H17").Value2 = 19.338559
$ws.Range("I17").Value2 = 0.1343775282351073
$ws.Range("J17").Value2 = 0.1343775282351073
$ws.Range("K17").Value2 = 1
$ws.Range("L17").Value2 = 0.3333333333333333
$ws.Range("M17").Value2 = 0.057141
$ws.Range("N17").Value2 = 0.171423
$ws.Range("O17").Value2 = 0.0004474116319810314
$ws.Range("P17").Value2 = 0.0004474116319810314
$ws.Range("Q17").Value2 = 0.368341533273
$ws.Range("R17").Value2 = 3.315073799457
$ws.Range("S17").Value2 = 0.00006012206920924647
$ws.Range("T17").Value2 = 0.00006012206920924649
